$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4471.4443
$ws.Range("J29").Value = 6491.25
$ws.Range("L29").Value = 19473.75
$ws.Range("N29").Value = -20035.75
$ws.Range("H87").Value = 33344.184
$ws.Range("J87").Value = 33344.184
$ws.Range("L87").Value = 33344.184
$ws.Range("N87").Value = -35840.184
$ws.Range("H90").Value = 33344.184
$ws.Range("J90").Value = 33344.184
$ws.Range("L90").Value = 100032.552
$ws.Range("N90").Value = -112512.552
$ws.Range("H100").Value = 3104.8125
$ws.Range("I100").Value = 2529.4
$ws.Range("J100").Value = 4063.8333
$ws.Range("K100").Value = 2529.4
$ws.Range("L100").Value = 4063.8333
$ws.Range("M100").Value = -1988.4
$ws.Range("N100").Value = -5145.8333
$ws.Range("H132").Value = 4662.304
$ws.Range("I132").Value = 1437.5625
$ws.Range("J132").Value = 12033.143
$ws.Range("K132").Value = 4312.6875
$ws.Range("L132").Value = 36099.429
$ws.Range("M132").Value = -1782.6875
$ws.Range("N132").Value = -41159.429
$ws.Range("H137").Value = 73280
$ws.Range("I137").Value = 2094.5715
$ws.Range("J137").Value = 144465.42
$ws.Range("K137").Value = 6283.7145
$ws.Range("L137").Value = 433396.26
$ws.Range("M137").Value = -3733.7145
$ws.Range("N137").Value = -438496.26
$ws.Range("H141").Value = 103104.11
$ws.Range("I141").Value = 103104.11
$ws.Range("K141").Value = 309312.33
$ws.Range("M141").Value = -304132.33

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 21935.375
$ws.Range("I45").Value = 22545.7
$ws.Range("J45").Value = 20918.166
$ws.Range("K45").Value = 22545.7
$ws.Range("L45").Value = 20918.166
$ws.Range("M45").Value = -22168.7
$ws.Range("N45").Value = -21672.166
$ws.Range("H63").Value = 1883.3334
$ws.Range("I63").Value = 1883.3334
$ws.Range("K63").Value = 1883.3334
$ws.Range("M63").Value = -1197.3334
$ws.Range("H66").Value = 1883.3334
$ws.Range("I66").Value = 1883.3334
$ws.Range("K66").Value = 9416.666999999999
$ws.Range("M66").Value = -5984.666999999999
$ws.Range("H97").Value = 2349.5
$ws.Range("I97").Value = 1971.6666
$ws.Range("K97").Value = 1971.6666
$ws.Range("M97").Value = -1475.6666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10510.2
$ws.Range("I20").Value = 8901.799999999999
$ws.Range("K20").Value = 8901.799999999999
$ws.Range("M20").Value = -8654.799999999999
$ws.Range("H99").Value = 3069.125
$ws.Range("I99").Value = 3069.125
$ws.Range("K99").Value = 3069.125
$ws.Range("M99").Value = -1571.125

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1002500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1002500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1002500
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -1002726
$ws.Range("H6").Value = 5278353
$ws.Range("I6").Value = 7037137.5
$ws.Range("J6").Value = 2000
$ws.Range("K6").Value = 7037137.5
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = -7037024.5
$ws.Range("N6").Value = -2226
$ws.Range("H31").Value = 3084
$ws.Range("I31").Value = 2814.818
$ws.Range("K31").Value = 2814.818
$ws.Range("M31").Value = -2519.818
$ws.Range("H34").Value = 3084
$ws.Range("I34").Value = 2814.818
$ws.Range("K34").Value = 2814.818
$ws.Range("M34").Value = -2612.818
$ws.Range("H50").Value = 8499.299999999999
$ws.Range("J50").Value = 8499.299999999999
$ws.Range("L50").Value = 8499.299999999999
$ws.Range("N50").Value = -9749.299999999999
$ws.Range("H51").Value = 8749.25
$ws.Range("J51").Value = 8749.25
$ws.Range("L51").Value = 8749.25
$ws.Range("N51").Value = -10221.25
$ws.Range("H59").Value = 16441.111
$ws.Range("J59").Value = 16871.25
$ws.Range("L59").Value = 16871.25
$ws.Range("N59").Value = -19161.25
$ws.Range("H60").Value = 8462.454
$ws.Range("J60").Value = 8332.666999999999
$ws.Range("L60").Value = 8332.666999999999
$ws.Range("N60").Value = -9354.666999999999
$ws.Range("H61").Value = 8749.25
$ws.Range("J61").Value = 8749.25
$ws.Range("L61").Value = 8749.25
$ws.Range("N61").Value = -9445.25
$ws.Range("H62").Value = 3099
$ws.Range("I62").Value = 3099
$ws.Range("K62").Value = 3099
$ws.Range("M62").Value = -2475
$ws.Range("H65").Value = 3099
$ws.Range("I65").Value = 3099
$ws.Range("K65").Value = 15495
$ws.Range("M65").Value = -12375
$ws.Range("H68").Value = 24996.363
$ws.Range("J68").Value = 24996.363
$ws.Range("L68").Value = 24996.363
$ws.Range("N68").Value = -26494.363
$ws.Range("H71").Value = 24996.363
$ws.Range("J71").Value = 24996.363
$ws.Range("L71").Value = 74989.08900000001
$ws.Range("N71").Value = -82477.08900000001
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = ""
$ws.Range("H103").Value = 23073.666
$ws.Range("I103").Value = 24110.5
$ws.Range("J103").Value = 21000
$ws.Range("K103").Value = 24110.5
$ws.Range("L103").Value = 21000
$ws.Range("M103").Value = -22938.5
$ws.Range("N103").Value = -23344

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 13120.5
$ws.Range("I62").Value = 11249.5
$ws.Range("J62").Value = 14991.5
$ws.Range("K62").Value = 33748.5
$ws.Range("L62").Value = 44974.5
$ws.Range("M62").Value = -33062.5
$ws.Range("N62").Value = -46346.5
$ws.Range("H65").Value = 13120.5
$ws.Range("I65").Value = 11249.5
$ws.Range("J65").Value = 14991.5
$ws.Range("K65").Value = 101245.5
$ws.Range("L65").Value = 134923.5
$ws.Range("M65").Value = -97813.5
$ws.Range("N65").Value = -141787.5
$ws.Range("H121").Value = 22224268
$ws.Range("J121").Value = 2914.3
$ws.Range("L121").Value = 8742.900000000001
$ws.Range("N121").Value = -11362.9
$ws.Range("H140").Value = 4357.731
$ws.Range("I140").Value = 2378.9443
$ws.Range("K140").Value = 7136.8329
$ws.Range("M140").Value = -1956.8329

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = ""
$ws.Range("H132").Value = 924.75
$ws.Range("I132").Value = 924.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2774.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -244.25
$ws.Range("N132").Value = ""

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 14349.643
$ws.Range("I122").Value = 15949.583
$ws.Range("J122").Value = 4750
$ws.Range("K122").Value = 47848.749
$ws.Range("L122").Value = 14250
$ws.Range("M122").Value = -45398.749
$ws.Range("N122").Value = -19150

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("H81").Value = 3658.5
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 3658.5
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H122").Value = 2253.2856
$ws.Range("I122").Value = 2253.2856
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6759.8568
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4309.8568
$ws.Range("N122").Value = ""
$ws.Range("H132").Value = 5715.5
$ws.Range("I132").Value = 7824.6875
$ws.Range("K132").Value = 23474.0625
$ws.Range("M132").Value = -20944.0625
